$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") rows 2 through 38 hold a date serial value that was
# bumped by one day (45696 -> 45697) during the automatic update.
for ($row = 2; $row -le 38; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45696) {
        $cell.Value2 = 45697
    }
}
